# Apply edits described in the diff:
# 1. Append 8 new daily data rows (112-119) to "coronadata_age_sex" sheet (sheet1)
# 2. Update summary totals on "gender breakdown" sheet (sheet2) rows 2-3
# 3. Update view/selection state: sheet1 selection moves to V122 with pane scrolled near new rows,
#    and the "gender breakdown" sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Append new rows 112-119 to sheet1 ---
$ws1.Range("A112").Value2 = 44019
$ws1.Range("B112").Value2 = 99
$ws1.Range("C112").Value2 = 131
$ws1.Range("D112").Value2 = 237
$ws1.Range("E112").Value2 = 442
$ws1.Range("F112").Value2 = 210
$ws1.Range("G112").Value2 = 357
$ws1.Range("H112").Value2 = 152
$ws1.Range("I112").Value2 = 287
$ws1.Range("J112").Value2 = 218
$ws1.Range("K112").Value2 = 308
$ws1.Range("L112").Value2 = 210
$ws1.Range("M112").Value2 = 237
$ws1.Range("N112").Value2 = 179
$ws1.Range("O112").Value2 = 202
$ws1.Range("P112").Value2 = 120
$ws1.Range("Q112").Value2 = 184
$ws1.Range("R112").Value2 = 59
$ws1.Range("S112").Value2 = 161
$ws1.Range("T112").Value2 = 3
$ws1.Range("U112").Value2 = 10
$ws1.Range("V112").Formula = "=SUM(B112,D112,F112,H112,J112,L112,N112,P112,R112,T112)"
$ws1.Range("W112").Formula = "=SUM(C112,E112,G112,I112,K112,M112,O112,Q112,S112,U112)"
$ws1.Range("X112").Formula = "=(V112/(V112+W112))*100"
$ws1.Range("Y112").Formula = "=(W112/(V112+W112))*100"

$ws1.Range("A113").Value2 = 44020
$ws1.Range("B113").Value2 = 101
$ws1.Range("C113").Value2 = 134
$ws1.Range("D113").Value2 = 241
$ws1.Range("E113").Value2 = 449
$ws1.Range("F113").Value2 = 213
$ws1.Range("G113").Value2 = 358
$ws1.Range("H113").Value2 = 153
$ws1.Range("I113").Value2 = 289
$ws1.Range("J113").Value2 = 220
$ws1.Range("K113").Value2 = 314
$ws1.Range("L113").Value2 = 210
$ws1.Range("M113").Value2 = 237
$ws1.Range("N113").Value2 = 180
$ws1.Range("O113").Value2 = 203
$ws1.Range("P113").Value2 = 120
$ws1.Range("Q113").Value2 = 184
$ws1.Range("R113").Value2 = 59
$ws1.Range("S113").Value2 = 161
$ws1.Range("T113").Value2 = 3
$ws1.Range("U113").Value2 = 10
$ws1.Range("V113").Formula = "=SUM(B113,D113,F113,H113,J113,L113,N113,P113,R113,T113)"
$ws1.Range("W113").Formula = "=SUM(C113,E113,G113,I113,K113,M113,O113,Q113,S113,U113)"
$ws1.Range("X113").Formula = "=(V113/(V113+W113))*100"
$ws1.Range("Y113").Formula = "=(W113/(V113+W113))*100"

$ws1.Range("A114").Value2 = 44021
$ws1.Range("B114").Value2 = 102
$ws1.Range("C114").Value2 = 134
$ws1.Range("D114").Value2 = 243
$ws1.Range("E114").Value2 = 465
$ws1.Range("F114").Value2 = 213
$ws1.Range("G114").Value2 = 360
$ws1.Range("H114").Value2 = 154
$ws1.Range("I114").Value2 = 292
$ws1.Range("J114").Value2 = 221
$ws1.Range("K114").Value2 = 316
$ws1.Range("L114").Value2 = 210
$ws1.Range("M114").Value2 = 238
$ws1.Range("N114").Value2 = 181
$ws1.Range("O114").Value2 = 204
$ws1.Range("P114").Value2 = 120
$ws1.Range("Q114").Value2 = 184
$ws1.Range("R114").Value2 = 59
$ws1.Range("S114").Value2 = 161
$ws1.Range("T114").Value2 = 3
$ws1.Range("U114").Value2 = 10
$ws1.Range("V114").Formula = "=SUM(B114,D114,F114,H114,J114,L114,N114,P114,R114,T114)"
$ws1.Range("W114").Formula = "=SUM(C114,E114,G114,I114,K114,M114,O114,Q114,S114,U114)"
$ws1.Range("X114").Formula = "=(V114/(V114+W114))*100"
$ws1.Range("Y114").Formula = "=(W114/(V114+W114))*100"

$ws1.Range("A115").Value2 = 44022
$ws1.Range("B115").Value2 = 102
$ws1.Range("C115").Value2 = 139
$ws1.Range("D115").Value2 = 248
$ws1.Range("E115").Value2 = 475
$ws1.Range("F115").Value2 = 222
$ws1.Range("G115").Value2 = 363
$ws1.Range("H115").Value2 = 155
$ws1.Range("I115").Value2 = 294
$ws1.Range("J115").Value2 = 227
$ws1.Range("K115").Value2 = 323
$ws1.Range("L115").Value2 = 213
$ws1.Range("M115").Value2 = 241
$ws1.Range("N115").Value2 = 181
$ws1.Range("O115").Value2 = 205
$ws1.Range("P115").Value2 = 120
$ws1.Range("Q115").Value2 = 184
$ws1.Range("R115").Value2 = 59
$ws1.Range("S115").Value2 = 161
$ws1.Range("T115").Value2 = 3
$ws1.Range("U115").Value2 = 10
$ws1.Range("V115").Formula = "=SUM(B115,D115,F115,H115,J115,L115,N115,P115,R115,T115)"
$ws1.Range("W115").Formula = "=SUM(C115,E115,G115,I115,K115,M115,O115,Q115,S115,U115)"
$ws1.Range("X115").Formula = "=(V115/(V115+W115))*100"
$ws1.Range("Y115").Formula = "=(W115/(V115+W115))*100"

$ws1.Range("A116").Value2 = 44023
$ws1.Range("B116").Value2 = 103
$ws1.Range("C116").Value2 = 140
$ws1.Range("D116").Value2 = 256
$ws1.Range("E116").Value2 = 480
$ws1.Range("F116").Value2 = 223
$ws1.Range("G116").Value2 = 367
$ws1.Range("H116").Value2 = 156
$ws1.Range("I116").Value2 = 297
$ws1.Range("J116").Value2 = 229
$ws1.Range("K116").Value2 = 325
$ws1.Range("L116").Value2 = 213
$ws1.Range("M116").Value2 = 242
$ws1.Range("N116").Value2 = 182
$ws1.Range("O116").Value2 = 206
$ws1.Range("P116").Value2 = 120
$ws1.Range("Q116").Value2 = 184
$ws1.Range("R116").Value2 = 59
$ws1.Range("S116").Value2 = 161
$ws1.Range("T116").Value2 = 3
$ws1.Range("U116").Value2 = 10
$ws1.Range("V116").Formula = "=SUM(B116,D116,F116,H116,J116,L116,N116,P116,R116,T116)"
$ws1.Range("W116").Formula = "=SUM(C116,E116,G116,I116,K116,M116,O116,Q116,S116,U116)"
$ws1.Range("X116").Formula = "=(V116/(V116+W116))*100"
$ws1.Range("Y116").Formula = "=(W116/(V116+W116))*100"

$ws1.Range("A117").Value2 = 44024
$ws1.Range("B117").Value2 = 107
$ws1.Range("C117").Value2 = 143
$ws1.Range("D117").Value2 = 258
$ws1.Range("E117").Value2 = 482
$ws1.Range("F117").Value2 = 226
$ws1.Range("G117").Value2 = 369
$ws1.Range("H117").Value2 = 157
$ws1.Range("I117").Value2 = 298
$ws1.Range("J117").Value2 = 230
$ws1.Range("K117").Value2 = 328
$ws1.Range("L117").Value2 = 215
$ws1.Range("M117").Value2 = 243
$ws1.Range("N117").Value2 = 183
$ws1.Range("O117").Value2 = 207
$ws1.Range("P117").Value2 = 120
$ws1.Range("Q117").Value2 = 185
$ws1.Range("R117").Value2 = 59
$ws1.Range("S117").Value2 = 161
$ws1.Range("T117").Value2 = 3
$ws1.Range("U117").Value2 = 10
$ws1.Range("V117").Formula = "=SUM(B117,D117,F117,H117,J117,L117,N117,P117,R117,T117)"
$ws1.Range("W117").Formula = "=SUM(C117,E117,G117,I117,K117,M117,O117,Q117,S117,U117)"
$ws1.Range("X117").Formula = "=(V117/(V117+W117))*100"
$ws1.Range("Y117").Formula = "=(W117/(V117+W117))*100"

$ws1.Range("A118").Value2 = 44025
$ws1.Range("B118").Value2 = 107
$ws1.Range("C118").Value2 = 144
$ws1.Range("D118").Value2 = 260
$ws1.Range("E118").Value2 = 485
$ws1.Range("F118").Value2 = 227
$ws1.Range("G118").Value2 = 371
$ws1.Range("H118").Value2 = 158
$ws1.Range("I118").Value2 = 302
$ws1.Range("J118").Value2 = 231
$ws1.Range("K118").Value2 = 329
$ws1.Range("L118").Value2 = 217
$ws1.Range("M118").Value2 = 245
$ws1.Range("N118").Value2 = 183
$ws1.Range("O118").Value2 = 207
$ws1.Range("P118").Value2 = 120
$ws1.Range("Q118").Value2 = 185
$ws1.Range("R118").Value2 = 59
$ws1.Range("S118").Value2 = 161
$ws1.Range("T118").Value2 = 3
$ws1.Range("U118").Value2 = 10
$ws1.Range("V118").Formula = "=SUM(B118,D118,F118,H118,J118,L118,N118,P118,R118,T118)"
$ws1.Range("W118").Formula = "=SUM(C118,E118,G118,I118,K118,M118,O118,Q118,S118,U118)"
$ws1.Range("X118").Formula = "=(V118/(V118+W118))*100"
$ws1.Range("Y118").Formula = "=(W118/(V118+W118))*100"

$ws1.Range("A119").Value2 = 44026
$ws1.Range("B119").Value2 = 108
$ws1.Range("C119").Value2 = 150
$ws1.Range("D119").Value2 = 270
$ws1.Range("E119").Value2 = 495
$ws1.Range("F119").Value2 = 228
$ws1.Range("G119").Value2 = 375
$ws1.Range("H119").Value2 = 161
$ws1.Range("I119").Value2 = 303
$ws1.Range("J119").Value2 = 232
$ws1.Range("K119").Value2 = 334
$ws1.Range("L119").Value2 = 218
$ws1.Range("M119").Value2 = 247
$ws1.Range("N119").Value2 = 185
$ws1.Range("O119").Value2 = 207
$ws1.Range("P119").Value2 = 120
$ws1.Range("Q119").Value2 = 185
$ws1.Range("R119").Value2 = 59
$ws1.Range("S119").Value2 = 161
$ws1.Range("T119").Value2 = 3
$ws1.Range("U119").Value2 = 10
$ws1.Range("V119").Formula = "=SUM(B119,D119,F119,H119,J119,L119,N119,P119,R119,T119)"
$ws1.Range("W119").Formula = "=SUM(C119,E119,G119,I119,K119,M119,O119,Q119,S119,U119)"
$ws1.Range("X119").Formula = "=(V119/(V119+W119))*100"
$ws1.Range("Y119").Formula = "=(W119/(V119+W119))*100"

# --- 2. Update gender breakdown totals on sheet2 ---
$ws2.Range("B2").Value2 = 1602
$ws2.Range("C2").Value2 = 379
$ws2.Range("D2").Value2 = 137
$ws2.Range("E2").Value2 = 122

$ws2.Range("B3").Value2 = 2583
$ws2.Range("C3").Value2 = 408
$ws2.Range("D3").Value2 = 100
$ws2.Range("E3").Value2 = 149

# --- 3. Update sheet1 view state: scroll frozen pane near the new rows and select V122 ---
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 96
$ws1.Range("V122").Select()

# --- 4. Activate "gender breakdown" sheet as the final active sheet/tab ---
$ws2.Activate()

Write-Host "done"
